# "use index culture date"
# Insert three new columns (O, P, Q) before the old "Polymicrobial_Infection"
# column so the remaining columns (old O:W) shift right to (new R:Z), then
# populate the three new columns with the Index_Culture + 72 hrs /
# Final_Result_Date / Final_Result - Index_Culture + 72 (hrs) data, and
# update the Cefepime / Total_DOT totals for row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank columns starting at O; this shifts old O:W -> R:Z and
# updates the sheet dimension automatically (A1:W2 -> A1:Z2).
$ws.Columns("O:Q").Insert()

# New header labels for the inserted columns.
$ws.Range("O1").Value = "Index_Culture + 72 hrs"
$ws.Range("P1").Value = "Final_Result_Date"
$ws.Range("Q1").Value = "Final_Result - Index_Culture + 72 (hrs)"

# New row 2 data for the inserted columns.
$ws.Range("N2").Value = 44581.62222222222
$ws.Range("O2").Value = 44584.62222222222
$ws.Range("P2").Value = 44585.58055555556
$ws.Range("Q2").Value = 23

# Match the date/time display format used elsewhere in the sheet (style
# used by Index_Culture's neighboring date cells, e.g. Last_Admin).
$ws.Range("N2:P2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Updated totals (columns shifted: old S2/V2 -> new V2/Y2).
$ws.Range("V2").Value = 10
$ws.Range("Y2").Value = 14
